$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "57.913.40"
$ws.Range("E2").Value = "  -3.90%  "
$ws.Range("D3").Value = "3.124.61"
$ws.Range("E3").Value = "  -5.43%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.06%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.120.82"
$ws.Range("E8").Value = "  -5.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.443"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.43%  "
$ws.Range("E10").Value = "  -7.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.109"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.383"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.01%  "
$ws.Range("D13").Value = "3.674.05"
$ws.Range("E13").Value = "  -5.01%  "
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.79%  "
$ws.Range("D16").Value = "3.122.24"
$ws.Range("E16").Value = "  -5.36%  "
$ws.Range("D17").Value = "57.808.55"
$ws.Range("E17").Value = "  -3.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000152"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "344.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.90%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.503"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.76%  "
$ws.Range("D26").Value = "3.257.36"
$ws.Range("E26").Value = "  -5.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.167"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("D28").Value = "0.0₃0952"
$ws.Range("E28").Value = "  -7.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.78"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.67%  "
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.97%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "21.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.64%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.70%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -10.01%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0691"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.01%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "24.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.82%  "
$ws.Range("D42").Value = "3.154.67"
$ws.Range("E42").Value = "  -5.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.695"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.47%  "
$ws.Range("E45").Value = "  -2.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.88%  "
$ws.Range("D49").Value = "2.254.04"
$ws.Range("E49").Value = "  -4.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.22%  "
